$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 545.75
$ws.Range("I31").Value = 545.75
$ws.Range("K31").Value = 1637.25
$ws.Range("M31").Value = -1407.25

$ws.Range("H38").Value = 1237.25
$ws.Range("I38").Value = 1237.25
$ws.Range("K38").Value = 3711.75
$ws.Range("M38").Value = -3339.75

$ws.Range("H39").Value = 407.0263
$ws.Range("I39").Value = 235.23077
$ws.Range("J39").Value = 496.36
$ws.Range("K39").Value = 705.69231
$ws.Range("L39").Value = 1489.08
$ws.Range("M39").Value = -409.69231
$ws.Range("N39").Value = -2081.08

$ws.Range("H51").Value = 78525.86
$ws.Range("J51").Value = 7333.625
$ws.Range("L51").Value = 7333.625
$ws.Range("N51").Value = -8301.625

$ws.Range("H61").Value = 315
$ws.Range("I61").Value = 315
$ws.Range("K61").Value = 945
$ws.Range("M61").Value = -773

$ws.Range("H97").Value = 1514.1666
$ws.Range("J97").Value = 1514.1666
$ws.Range("L97").Value = 4542.4998
$ws.Range("N97").Value = -5534.4998

$ws.Range("H111").Value = 5450
$ws.Range("I111").Value = 5900
$ws.Range("J111").Value = 5000
$ws.Range("K111").Value = 17700
$ws.Range("L111").Value = 15000
$ws.Range("M111").Value = -14633
$ws.Range("N111").Value = -21134

$ws.Range("H112").Value = 1186.4
$ws.Range("I112").Value = 474.875
$ws.Range("J112").Value = 1999.5714
$ws.Range("K112").Value = 1424.625
$ws.Range("L112").Value = 5998.7142
$ws.Range("M112").Value = -316.625
$ws.Range("N112").Value = -8214.7142

$ws.Range("H138").Value = 2897.6758
$ws.Range("I138").Value = 1032.7
$ws.Range("J138").Value = 3189.0781
$ws.Range("K138").Value = 3098.1
$ws.Range("L138").Value = 9567.2343
$ws.Range("M138").Value = 2041.9
$ws.Range("N138").Value = -19847.2343

$ws.Range("H141").Value = 2674.55
$ws.Range("I141").Value = 2696.8823
$ws.Range("J141").Value = 2548
$ws.Range("K141").Value = 8090.646900000001
$ws.Range("L141").Value = 7644
$ws.Range("M141").Value = -2910.646900000001
$ws.Range("N141").Value = -18004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 732.75
$ws.Range("J4").Value = 1289
$ws.Range("L4").Value = 1289
$ws.Range("N4").Value = -1521

$ws.Range("H32").Value = 10206482
$ws.Range("I32").Value = 11112731
$ws.Range("K32").Value = 11112731
$ws.Range("M32").Value = -11112444

$ws.Range("H133").Value = 69997.5
$ws.Range("J133").Value = 69997.5
$ws.Range("L133").Value = 69997.5
$ws.Range("N133").Value = -75057.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3714.8333
$ws.Range("I86").Value = 3998.75
$ws.Range("J86").Value = 3147
$ws.Range("K86").Value = 3998.75
$ws.Range("L86").Value = 3147
$ws.Range("M86").Value = -2875.75
$ws.Range("N86").Value = -5393

$ws.Range("H89").Value = 3714.8333
$ws.Range("I89").Value = 3998.75
$ws.Range("J89").Value = 3147
$ws.Range("K89").Value = 19993.75
$ws.Range("L89").Value = 15735
$ws.Range("M89").Value = -14377.75
$ws.Range("N89").Value = -26967

$ws.Range("H105").Value = 1932.2222
$ws.Range("I105").Value = 1078
$ws.Range("K105").Value = 1078
$ws.Range("M105").Value = 669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 502
$ws.Range("J86").Value = 602.6667
$ws.Range("L86").Value = 1808.0001
$ws.Range("N86").Value = -4180.0001

$ws.Range("H89").Value = 502
$ws.Range("J89").Value = 602.6667
$ws.Range("L89").Value = 5424.0003
$ws.Range("N89").Value = -17280.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 481.5
$ws.Range("I2").Value = 81.25
$ws.Range("J2").Value = 1282
$ws.Range("K2").Value = 81.25
$ws.Range("L2").Value = 1282
$ws.Range("M2").Value = 31.75
$ws.Range("N2").Value = -1508

$ws.Range("H70").Value = 5399.3335
$ws.Range("J70").Value = 5250
$ws.Range("L70").Value = 5250
$ws.Range("N70").Value = -5790

$ws.Range("H73").Value = 5399.3335
$ws.Range("J73").Value = 5250
$ws.Range("L73").Value = 5250
$ws.Range("N73").Value = -7122

$ws.Range("H97").Value = 1702.9412
$ws.Range("I97").Value = 1860.8572
$ws.Range("J97").Value = 966
$ws.Range("K97").Value = 1860.8572
$ws.Range("L97").Value = 966
$ws.Range("M97").Value = -1364.8572
$ws.Range("N97").Value = -1958

$ws.Range("H113").Value = 3282.9167
$ws.Range("I113").Value = 1837.5
$ws.Range("K113").Value = 1837.5
$ws.Range("M113").Value = 332.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2551228.5
$ws.Range("I43").Value = 3367000
$ws.Range("J43").Value = 1939400
$ws.Range("K43").Value = 3367000
$ws.Range("L43").Value = 1939400
$ws.Range("M43").Value = -3366807
$ws.Range("N43").Value = -1939786

$ws.Range("H93").Value = 142860450
$ws.Range("I93").Value = 200002820
$ws.Range("K93").Value = 200002820
$ws.Range("M93").Value = -200001572

$ws.Range("H100").Value = 3083.375
$ws.Range("I100").Value = 3213.8
$ws.Range("J100").Value = 2866
$ws.Range("K100").Value = 3213.8
$ws.Range("L100").Value = 2866
$ws.Range("M100").Value = -2672.8
$ws.Range("N100").Value = -3948

$ws.Range("H122").Value = 6070.294
$ws.Range("I122").Value = 5274.1665
$ws.Range("J122").Value = 7981
$ws.Range("K122").Value = 15822.4995
$ws.Range("L122").Value = 23943
$ws.Range("M122").Value = -13372.4995
$ws.Range("N122").Value = -28843

$ws.Range("H132").Value = 119939.586
$ws.Range("I132").Value = 64935.5
$ws.Range("K132").Value = 194806.5
$ws.Range("M132").Value = -192276.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H95").Value = 80000
$ws.Range("I95").Value = 60000
$ws.Range("J95").Value = 100000
$ws.Range("K95").Value = 60000
$ws.Range("L95").Value = 100000
$ws.Range("M95").Value = -57254
$ws.Range("N95").Value = -105492
